$wb = $excel.ActiveWorkbook
$wsTracking = $wb.Worksheets.Item("SprintTracking")
$wsBacklog = $wb.Worksheets.Item("Backlog")

# Add new backlog task descriptions (rows 2-8 of the Backlog sheet)
$wsBacklog.Range("A2").Value = "Authentication for users"
$wsBacklog.Range("A3").Value = "GUID to help with sensitive information"
$wsBacklog.Range("A4").Value = "Errors - error codes"
$wsBacklog.Range("A5").Value = "Root endpoint - root to GET all endpoint categories that the rest api supports"
$wsBacklog.Range("A6").Value = "HEAD - issued agains any endpoint to get just the header info."
$wsBacklog.Range("A7").Value = "Pagination - ?page, limit results returned to prevent dos attacks."
$wsBacklog.Range("A8").Value = "Schema - how to develop?"

# Update selections to match the author's last cursor position in each sheet
$wsBacklog.Range("A9").Select()
$wsTracking.Range("D15").Select()
